$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1861111111111111
$ws.Range("C2").Value = 0.5722222222222222
$ws.Range("J2").Value = 0.005555555555555556
$ws.Range("P2").Value = 0.1638888888888889
$ws.Range("S2").Value = 0.07222222222222222
$ws.Range("B3").Value = 0.02336448598130841
$ws.Range("C3").Value = 0.02336448598130841
$ws.Range("J3").Value = 0.009345794392523364
$ws.Range("P3").Value = 0.7663551401869159
$ws.Range("S3").Value = 0.1775700934579439
$ws.Range("J4").Value = 0.05084745762711865
$ws.Range("O4").Value = 0.01694915254237288
$ws.Range("P4").Value = 0.6610169491525424
$ws.Range("S4").Value = 0.2711864406779661
$ws.Range("B6").Value = 0.09259259259259259
$ws.Range("D6").Value = 0.01234567901234568
$ws.Range("F6").Value = 0.04938271604938271
$ws.Range("J6").Value = 0.2469135802469136
$ws.Range("O6").Value = 0.04938271604938271
$ws.Range("Q6").Value = 0.1419753086419753
$ws.Range("R6").Value = 0.06172839506172839
$ws.Range("S6").Value = 0.345679012345679
$ws.Range("B7").Value = 0.1764705882352941
$ws.Range("D7").Value = 0.05294117647058823
$ws.Range("F7").Value = 0.03529411764705882
$ws.Range("J7").Value = 0.1
$ws.Range("O7").Value = 0.01176470588235294
$ws.Range("Q7").Value = 0.1647058823529412
$ws.Range("R7").Value = 0.09411764705882353
$ws.Range("S7").Value = 0.3647058823529412
$ws.Range("B8").Value = 0.1192893401015228
$ws.Range("D8").Value = 0.03045685279187817
$ws.Range("F8").Value = 0.04060913705583756
$ws.Range("J8").Value = 0.07360406091370558
$ws.Range("O8").Value = 0.01269035532994924
$ws.Range("Q8").Value = 0.1928934010152284
$ws.Range("R8").Value = 0.09390862944162437
$ws.Range("S8").Value = 0.4365482233502538
$ws.Range("B9").Value = 0.1272727272727273
$ws.Range("D9").Value = 0.01818181818181818
$ws.Range("E9").Value = 0.004545454545454545
$ws.Range("F9").Value = 0.08181818181818182
$ws.Range("J9").Value = 0.06818181818181818
$ws.Range("Q9").Value = 0.1818181818181818
$ws.Range("R9").Value = 0.1181818181818182
$ws.Range("S9").Value = 0.4
$ws.Range("B10").Value = 0.1389645776566757
$ws.Range("D10").Value = 0.02997275204359673
$ws.Range("E10").Value = 0.0009082652134423251
$ws.Range("F10").Value = 0.05177111716621254
$ws.Range("J10").Value = 0.07447774750227067
$ws.Range("O10").Value = 0.0145322434150772
$ws.Range("Q10").Value = 0.2297910990009083
$ws.Range("R10").Value = 0.08537693006357856
$ws.Range("S10").Value = 0.3742052679382379
$ws.Range("G11").Value = 0.1464285714285714
$ws.Range("J11").Value = 0.08928571428571429
$ws.Range("K11").Value = 0.2071428571428572
$ws.Range("L11").Value = 0.5428571428571428
$ws.Range("S11").Value = 0.01428571428571429
$ws.Range("G12").Value = 0.7088607594936709
$ws.Range("J12").Value = 0.2025316455696203
$ws.Range("K12").Value = 0.03164556962025317
$ws.Range("L12").Value = 0.0379746835443038
$ws.Range("S12").Value = 0.0189873417721519
$ws.Range("G13").Value = 0.5675675675675675
$ws.Range("J13").Value = 0.4054054054054054
$ws.Range("S13").Value = 0.02702702702702703
$ws.Range("F15").Value = 0.01657458563535912
$ws.Range("H15").Value = 0.1767955801104972
$ws.Range("I15").Value = 0.1215469613259668
$ws.Range("J15").Value = 0.3535911602209945
$ws.Range("K15").Value = 0.03314917127071823
$ws.Range("M15").Value = 0.005524861878453038
$ws.Range("O15").Value = 0.02209944751381215
$ws.Range("S15").Value = 0.2707182320441989
$ws.Range("F16").Value = 0.01968503937007874
$ws.Range("H16").Value = 0.1535433070866142
$ws.Range("I16").Value = 0.06299212598425197
$ws.Range("J16").Value = 0.4330708661417323
$ws.Range("K16").Value = 0.1102362204724409
$ws.Range("M16").Value = 0.007874015748031496
$ws.Range("N16").Value = 0.003937007874015748
$ws.Range("O16").Value = 0.04724409448818898
$ws.Range("S16").Value = 0.1614173228346457
$ws.Range("F17").Value = 0.0218978102189781
$ws.Range("H17").Value = 0.1508515815085158
$ws.Range("I17").Value = 0.1167883211678832
$ws.Range("J17").Value = 0.4501216545012166
$ws.Range("K17").Value = 0.0851581508515815
$ws.Range("M17").Value = 0.009732360097323601
$ws.Range("O17").Value = 0.0462287104622871
$ws.Range("S17").Value = 0.1192214111922141
$ws.Range("F18").Value = 0.03867403314917127
$ws.Range("H18").Value = 0.1491712707182321
$ws.Range("I18").Value = 0.1160220994475138
$ws.Range("J18").Value = 0.4088397790055249
$ws.Range("K18").Value = 0.09944751381215469
$ws.Range("M18").Value = 0.01104972375690608
$ws.Range("N18").Value = 0.005524861878453038
$ws.Range("O18").Value = 0.05524861878453038
$ws.Range("S18").Value = 0.1160220994475138
$ws.Range("F19").Value = 0.01005867560771165
$ws.Range("H19").Value = 0.1978206202849958
$ws.Range("I19").Value = 0.09388097233864208
$ws.Range("J19").Value = 0.3495389773679799
$ws.Range("K19").Value = 0.1081307627829002
$ws.Range("M19").Value = 0.02598491198658843
$ws.Range("N19").Value = 0.0008382229673093043
$ws.Range("O19").Value = 0.06789606035205364
$ws.Range("S19").Value = 0.1458507963118189
